# Adds columns I ("I0") and J ("IF") to the sheet, matching the style of
# the existing header row and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, matching the style used by the other header cells (B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row values for columns I and J (rows 2..73).
$iValues = @(9,9,9,8,8,10,9,9,9,9,9,10,9,9,9,10,8,9,8,9,8,9,8,9,9,8,8,8,8,7,7,8,8,8,9,9,8,9,8,9,9,8,8,8,9,10,8,8,7,8,8,7,7,9,8,7,8,8,9,9,8,8,9,7,7,7,8,7,9,7,8,6)
$jValues = @(9,9,9,8,8,10,9,9,9,9,9,10,9,9,9,10,8,9,8,9,8,9,9,9,9,8,9,8,8,7,7,8,8,8,9,9,8,9,8,9,9,8,8,8,9,10,8,8,7,8,8,7,7,9,8,7,8,8,9,9,8,8,9,7,7,7,8,7,9,7,8,6)

for ($r = 2; $r -le 73; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
